$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 55, pushing the existing rows 55:75 down to 56:76
$ws.Rows.Item(55).Insert()

# Populate the newly inserted row 55 with the new weekly price record
$ws.Cells.Item(55,1).Value2  = 4
$ws.Cells.Item(55,2).Value2  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(55,3).Value2  = "Los Lagos"
$ws.Cells.Item(55,4).Value2  = 44489
$ws.Cells.Item(55,5).Value2  = 10
$ws.Cells.Item(55,6).Value2  = 100112022
$ws.Cells.Item(55,7).Value2  = "Arveja Verde"
$ws.Cells.Item(55,8).Value2  = "Sin especificar"
$ws.Cells.Item(55,9).Value2  = "Primera"
$ws.Cells.Item(55,10).Value2 = 35
$ws.Cells.Item(55,11).Value2 = 25000
$ws.Cells.Item(55,12).Value2 = 25000
$ws.Cells.Item(55,13).Value2 = 25000
$ws.Cells.Item(55,14).Value2 = "`$/saco 25 kilos"
$ws.Cells.Item(55,15).Value2 = "Región del Maule"
$ws.Cells.Item(55,16).Value2 = 1000
$ws.Cells.Item(55,17).Value2 = 25
$ws.Cells.Item(55,18).Value2 = "Hortaliza"

# Keep the date column formatted the same way as the rest of column D
$ws.Cells.Item(55,4).NumberFormat = $ws.Cells.Item(56,4).NumberFormat
